# Updates the cryptos price/volume table to the latest scraped snapshot.
# Mirrors the GitHub Actions commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.299.08"
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").Value = "1.869.97"
$ws.Range("E3").Value = "  +0.47%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.06"
$ws.Range("E5").Value = "  -0.56%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4701"
$ws.Range("E7").Value = "  +0.46%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2870"
$ws.Range("E8").Value = "  +0.35%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06602"
$ws.Range("E9").Value = "  +1.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.63"
$ws.Range("E10").Value = "  -1.14%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07960"
$ws.Range("E11").Value = "  +0.52%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "96.77"
$ws.Range("E12").Value = "  -0.57%  "
$ws.Range("D13").Value = "1.876.23"
$ws.Range("E13").Value = "  +0.72%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6990"
$ws.Range("E14").Value = "  +2.57%  "
$ws.Range("E15").Value = "  -0.93%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "268.83"
$ws.Range("E16").Value = "  -0.60%  "
$ws.Range("D17").Value = "30.370.44"
$ws.Range("E17").Value = "  +0.50%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.13"
$ws.Range("E18").Value = "  +4.42%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007800"
$ws.Range("E19").Value = "  +6.16%  "
$ws.Range("E20").Value = "  -0.08%  "
$ws.Range("D21").Value = "2.123.88"
$ws.Range("E21").Value = "  +0.49%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.280"
$ws.Range("E23").Value = "  -0.88%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.220"
$ws.Range("E24").Value = "  +0.89%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.384"
$ws.Range("E25").Value = "  +1.85%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "167.50"
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.88"
$ws.Range("E27").Value = "  -0.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.953"
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("E29").Value = "  -1.38%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09919"
$ws.Range("E30").Value = "  +0.60%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.336"
$ws.Range("E31").Value = "  -0.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.460"
$ws.Range("E32").Value = "  -1.25%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.060"
$ws.Range("E33").Value = "  +0.13%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04714"
$ws.Range("E34").Value = "  +0.19%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7035"
$ws.Range("E36").Value = "  +0.45%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.727"
$ws.Range("E37").Value = "  +0.57%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01876"
$ws.Range("E38").Value = "  +0.35%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.808"
$ws.Range("E39").Value = "  +7.00%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.266"
$ws.Range("E40").Value = "  +0.06%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "72.17"
$ws.Range("E41").Value = "  -4.36%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.962"
$ws.Range("E42").Value = "  +1.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4183"
$ws.Range("E43").Value = "  +0.62%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8416"
$ws.Range("E44").Value = "  -1.20%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9994"
$ws.Range("E45").Value = "  -0.14%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "102.83"
$ws.Range("E46").Value = "  -0.38%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.117"
$ws.Range("E47").Value = "  -0.61%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.174"
$ws.Range("E48").Value = "  -0.93%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "919.79"
$ws.Range("E49").Value = "  -3.71%  "
$ws.Range("E50").Value = "  +1.46%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05681"
$ws.Range("E51").Value = "  +0.56%  "
